$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H18").Value = 500
$wsALC.Range("I18").Value = 500
$wsALC.Range("J18").Value = 0
$wsALC.Range("K18").Value = 500
$wsALC.Range("L18").Value = 0
$wsALC.Range("M18").Value = -216

$wsALC.Range("H40").Value = 4820
$wsALC.Range("I40").Value = 4700
$wsALC.Range("J40").Value = 4838.4614
$wsALC.Range("K40").Value = 4700
$wsALC.Range("L40").Value = 4838.4614
$wsALC.Range("M40").Value = -4525
$wsALC.Range("N40").Value = -5188.4614

$wsALC.Range("H41").Value = 1618.1666
$wsALC.Range("I41").Value = 1328
$wsALC.Range("J41").Value = 2198.5
$wsALC.Range("K41").Value = 1328
$wsALC.Range("L41").Value = 2198.5
$wsALC.Range("M41").Value = -888
$wsALC.Range("N41").Value = -3078.5

$wsALC.Range("H53").Value = 425.3158
$wsALC.Range("I53").Value = 404.64706
$wsALC.Range("J53").Value = 601
$wsALC.Range("K53").Value = 404.64706
$wsALC.Range("L53").Value = 601
$wsALC.Range("M53").Value = 232.35294

$wsALC.Range("H58").Value = 5400
$wsALC.Range("I58").Value = 250
$wsALC.Range("J58").Value = 8833.333000000001
$wsALC.Range("K58").Value = 750
$wsALC.Range("L58").Value = 26499.999
$wsALC.Range("M58").Value = -600
$wsALC.Range("N58").Value = -26799.999

$wsALC.Range("H70").Value = 3741.5833
$wsALC.Range("I70").Value = 5241.5
$wsALC.Range("J70").Value = 2241.6667
$wsALC.Range("K70").Value = 15724.5
$wsALC.Range("L70").Value = 6725.000100000001
$wsALC.Range("M70").Value = -15454.5
$wsALC.Range("N70").Value = -7265.000100000001

$wsALC.Range("H73").Value = 3741.5833
$wsALC.Range("I73").Value = 5241.5
$wsALC.Range("J73").Value = 2241.6667
$wsALC.Range("K73").Value = 15724.5
$wsALC.Range("L73").Value = 6725.000100000001
$wsALC.Range("M73").Value = -14788.5
$wsALC.Range("N73").Value = -8597.000100000001

$wsALC.Range("H96").Value = 175
$wsALC.Range("I96").Value = 200
$wsALC.Range("J96").Value = 150
$wsALC.Range("K96").Value = 600
$wsALC.Range("L96").Value = 450
$wsALC.Range("M96").Value = 773
$wsALC.Range("N96").Value = -3196

$wsALC.Range("H107").Value = 1039.6765
$wsALC.Range("I107").Value = 755.4286
$wsALC.Range("J107").Value = 2366.1667
$wsALC.Range("K107").Value = 755.4286
$wsALC.Range("L107").Value = 2366.1667
$wsALC.Range("M107").Value = 1164.5714
$wsALC.Range("N107").Value = -6206.1667

$wsALC.Range("H112").Value = 3068.8
$wsALC.Range("I112").Value = 1949.5
$wsALC.Range("J112").Value = 3348.625
$wsALC.Range("K112").Value = 5848.5
$wsALC.Range("L112").Value = 10045.875
$wsALC.Range("M112").Value = -4740.5
$wsALC.Range("N112").Value = -12261.875

$wsALC.Range("H132").Value = 3916.9714
$wsALC.Range("I132").Value = 3520.258
$wsALC.Range("J132").Value = 6991.5
$wsALC.Range("K132").Value = 10560.774
$wsALC.Range("L132").Value = 20974.5
$wsALC.Range("M132").Value = -8030.773999999999

$wsALC.Range("H133").Value = 300000
$wsALC.Range("I133").Value = 0
$wsALC.Range("J133").Value = 300000
$wsALC.Range("K133").Value = 0
$wsALC.Range("L133").Value = 300000
$wsALC.Range("N133").Value = -310120

$wsALC.Range("H138").Value = 3977.1365
$wsALC.Range("I138").Value = 2649.3333
$wsALC.Range("J138").Value = 4475.0625
$wsALC.Range("K138").Value = 7947.999899999999
$wsALC.Range("L138").Value = 13425.1875
$wsALC.Range("M138").Value = -2807.999899999999
$wsALC.Range("N138").Value = -23705.1875

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H31").Value = 6473.3335
$wsARM.Range("I31").Value = 6473.3335
$wsARM.Range("J31").Value = 0
$wsARM.Range("K31").Value = 6473.3335
$wsARM.Range("L31").Value = 0
$wsARM.Range("M31").Value = -6179.3335

$wsARM.Range("H32").Value = 25923
$wsARM.Range("I32").Value = 18180.615
$wsARM.Range("J32").Value = 76248.5
$wsARM.Range("K32").Value = 18180.615
$wsARM.Range("L32").Value = 76248.5
$wsARM.Range("M32").Value = -17893.615
$wsARM.Range("N32").Value = -76822.5

$wsARM.Range("H45").Value = 1456
$wsARM.Range("I45").Value = 1456
$wsARM.Range("J45").Value = 0
$wsARM.Range("K45").Value = 1456
$wsARM.Range("L45").Value = 0
$wsARM.Range("M45").Value = -1079

$wsARM.Range("H97").Value = 1469.8334
$wsARM.Range("I97").Value = 1469.8334
$wsARM.Range("J97").Value = 0
$wsARM.Range("K97").Value = 1469.8334
$wsARM.Range("L97").Value = 0
$wsARM.Range("M97").Value = -973.8334

$wsARM.Range("H110").Value = 4085.1667
$wsARM.Range("I110").Value = 2011
$wsARM.Range("J110").Value = 4500
$wsARM.Range("K110").Value = 2011
$wsARM.Range("L110").Value = 4500
$wsARM.Range("M110").Value = 34
$wsARM.Range("N110").Value = -8590

$wsARM.Range("H128").Value = 0
$wsARM.Range("I128").Value = 0
$wsARM.Range("J128").Value = 0
$wsARM.Range("K128").Value = 0
$wsARM.Range("L128").Value = 0
$wsARM.Range("N128").ClearContents()

$wsARM.Range("H139").Value = 45000
$wsARM.Range("I139").Value = 0
$wsARM.Range("J139").Value = 45000
$wsARM.Range("K139").Value = 0
$wsARM.Range("L139").Value = 45000
$wsARM.Range("N139").Value = -55280

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H81").Value = 55000
$wsBSM.Range("I81").Value = 0
$wsBSM.Range("J81").Value = 55000
$wsBSM.Range("K81").Value = 0
$wsBSM.Range("L81").Value = 55000
$wsBSM.Range("N81").Value = -57122

$wsBSM.Range("H84").Value = 55000
$wsBSM.Range("I84").Value = 0
$wsBSM.Range("J84").Value = 55000
$wsBSM.Range("K84").Value = 0
$wsBSM.Range("L84").Value = 165000
$wsBSM.Range("N84").Value = -175608

$wsBSM.Range("H86").Value = 0
$wsBSM.Range("I86").Value = 0
$wsBSM.Range("J86").Value = 0
$wsBSM.Range("K86").Value = 0
$wsBSM.Range("L86").Value = 0
$wsBSM.Range("M86").ClearContents()
$wsBSM.Range("N86").ClearContents()

$wsBSM.Range("H89").Value = 0
$wsBSM.Range("I89").Value = 0
$wsBSM.Range("J89").Value = 0
$wsBSM.Range("K89").Value = 0
$wsBSM.Range("L89").Value = 0
$wsBSM.Range("M89").ClearContents()
$wsBSM.Range("N89").ClearContents()

$wsBSM.Range("H94").Value = 1436.7693
$wsBSM.Range("I94").Value = 1294.24
$wsBSM.Range("J94").Value = 5000
$wsBSM.Range("K94").Value = 1294.24
$wsBSM.Range("L94").Value = 5000
$wsBSM.Range("M94").Value = -843.24

$wsBSM.Range("H107").Value = 2035.2858
$wsBSM.Range("I107").Value = 949.4
$wsBSM.Range("J107").Value = 4750
$wsBSM.Range("K107").Value = 949.4
$wsBSM.Range("L107").Value = 4750
$wsBSM.Range("M107").Value = 970.6

$wsBSM.Range("H135").Value = 0
$wsBSM.Range("I135").Value = 0
$wsBSM.Range("J135").Value = 0
$wsBSM.Range("K135").Value = 0
$wsBSM.Range("L135").Value = 0
$wsBSM.Range("N135").ClearContents()

$wsBSM.Range("H138").Value = 49998
$wsBSM.Range("I138").Value = 0
$wsBSM.Range("J138").Value = 49998
$wsBSM.Range("K138").Value = 0
$wsBSM.Range("L138").Value = 49998
$wsBSM.Range("N138").Value = -60278

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H31").Value = 7078.5
$wsCRP.Range("I31").Value = 5999
$wsCRP.Range("J31").Value = 7438.3335
$wsCRP.Range("K31").Value = 5999
$wsCRP.Range("L31").Value = 7438.3335
$wsCRP.Range("M31").Value = -5704

$wsCRP.Range("H34").Value = 7078.5
$wsCRP.Range("I34").Value = 5999
$wsCRP.Range("J34").Value = 7438.3335
$wsCRP.Range("K34").Value = 5999
$wsCRP.Range("L34").Value = 7438.3335
$wsCRP.Range("M34").Value = -5797

$wsCRP.Range("H44").Value = 23354.334
$wsCRP.Range("I44").Value = 23354.334
$wsCRP.Range("J44").Value = 0
$wsCRP.Range("K44").Value = 23354.334
$wsCRP.Range("L44").Value = 0
$wsCRP.Range("M44").Value = -22912.334

$wsCRP.Range("H58").Value = 6240.5264
$wsCRP.Range("I58").Value = 5149.8237
$wsCRP.Range("J58").Value = 15511.5
$wsCRP.Range("K58").Value = 5149.8237
$wsCRP.Range("L58").Value = 15511.5
$wsCRP.Range("M58").Value = -4946.8237

$wsCRP.Range("H86").Value = 3181.7646
$wsCRP.Range("I86").Value = 3131.25
$wsCRP.Range("J86").Value = 3990
$wsCRP.Range("K86").Value = 3131.25
$wsCRP.Range("L86").Value = 3990
$wsCRP.Range("M86").Value = -2008.25

$wsCRP.Range("H89").Value = 3181.7646
$wsCRP.Range("I89").Value = 3131.25
$wsCRP.Range("J89").Value = 3990
$wsCRP.Range("K89").Value = 15656.25
$wsCRP.Range("L89").Value = 19950
$wsCRP.Range("M89").Value = -10040.25

$wsCRP.Range("H94").Value = 1797.8334
$wsCRP.Range("I94").Value = 1998.6666
$wsCRP.Range("J94").Value = 1597
$wsCRP.Range("K94").Value = 1998.6666
$wsCRP.Range("L94").Value = 1597
$wsCRP.Range("M94").Value = -1547.6666
$wsCRP.Range("N94").Value = -2499

$wsCRP.Range("H112").Value = 65000
$wsCRP.Range("I112").Value = 0
$wsCRP.Range("J112").Value = 65000
$wsCRP.Range("K112").Value = 0
$wsCRP.Range("L112").Value = 65000
$wsCRP.Range("N112").Value = -67954

$wsCRP.Range("H121").Value = 40000
$wsCRP.Range("I121").Value = 0
$wsCRP.Range("J121").Value = 40000
$wsCRP.Range("K121").Value = 0
$wsCRP.Range("L121").Value = 40000
$wsCRP.Range("N121").Value = -42620
$wsCRP.Range("M121").ClearContents()

$wsCRP.Range("H135").Value = 54000
$wsCRP.Range("I135").Value = 0
$wsCRP.Range("J135").Value = 54000
$wsCRP.Range("K135").Value = 0
$wsCRP.Range("L135").Value = 54000
$wsCRP.Range("N135").Value = -64140

$wsCRP.Range("H136").Value = 6240.5264
$wsCRP.Range("I136").Value = 5149.8237
$wsCRP.Range("J136").Value = 15511.5
$wsCRP.Range("K136").Value = 15449.4711
$wsCRP.Range("L136").Value = 46534.5
$wsCRP.Range("M136").Value = -12899.4711

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H58").Value = 0
$wsCUL.Range("I58").Value = 0
$wsCUL.Range("J58").Value = 0
$wsCUL.Range("K58").Value = 0
$wsCUL.Range("L58").Value = 0
$wsCUL.Range("N58").ClearContents()

$wsCUL.Range("H129").Value = 3936.889
$wsCUL.Range("I129").Value = 0
$wsCUL.Range("J129").Value = 3936.889
$wsCUL.Range("K129").Value = 0
$wsCUL.Range("L129").Value = 11810.667
$wsCUL.Range("N129").Value = -21810.667

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H39").Value = 0
$wsGSM.Range("I39").Value = 0
$wsGSM.Range("J39").Value = 0
$wsGSM.Range("K39").Value = 0
$wsGSM.Range("L39").Value = 0
$wsGSM.Range("N39").ClearContents()

$wsGSM.Range("H69").Value = 106665
$wsGSM.Range("I69").Value = 0
$wsGSM.Range("J69").Value = 106665
$wsGSM.Range("K69").Value = 0
$wsGSM.Range("L69").Value = 106665
$wsGSM.Range("N69").Value = -108163

$wsGSM.Range("H72").Value = 106665
$wsGSM.Range("I72").Value = 0
$wsGSM.Range("J72").Value = 106665
$wsGSM.Range("K72").Value = 0
$wsGSM.Range("L72").Value = 319995
$wsGSM.Range("N72").Value = -327483

$wsGSM.Range("H80").Value = 4132.5
$wsGSM.Range("I80").Value = 0
$wsGSM.Range("J80").Value = 4132.5
$wsGSM.Range("K80").Value = 0
$wsGSM.Range("L80").Value = 4132.5
$wsGSM.Range("N80").Value = -6128.5

$wsGSM.Range("H83").Value = 4132.5
$wsGSM.Range("I83").Value = 0
$wsGSM.Range("J83").Value = 4132.5
$wsGSM.Range("K83").Value = 0
$wsGSM.Range("L83").Value = 20662.5
$wsGSM.Range("N83").Value = -30646.5

$wsGSM.Range("H123").Value = 5024999
$wsGSM.Range("I123").Value = 0
$wsGSM.Range("J123").Value = 5024999
$wsGSM.Range("K123").Value = 0
$wsGSM.Range("L123").Value = 5024999
$wsGSM.Range("N123").Value = -5029899

$wsGSM.Range("H129").Value = 50000
$wsGSM.Range("I129").Value = 50000
$wsGSM.Range("J129").Value = 0
$wsGSM.Range("K129").Value = 50000
$wsGSM.Range("L129").Value = 0
$wsGSM.Range("M129").Value = -45000

$wsGSM.Range("H133").Value = 75321
$wsGSM.Range("I133").Value = 0
$wsGSM.Range("J133").Value = 75321
$wsGSM.Range("K133").Value = 0
$wsGSM.Range("L133").Value = 75321
$wsGSM.Range("N133").Value = -85441

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H40").Value = 4341.9375
$wsLTW.Range("I40").Value = 4544.357
$wsLTW.Range("J40").Value = 2925
$wsLTW.Range("K40").Value = 4544.357
$wsLTW.Range("L40").Value = 2925
$wsLTW.Range("M40").Value = -4408.357
$wsLTW.Range("N40").Value = -3197

$wsLTW.Range("H132").Value = 13037.277
$wsLTW.Range("I132").Value = 14690.857
$wsLTW.Range("J132").Value = 7249.75
$wsLTW.Range("K132").Value = 44072.571
$wsLTW.Range("L132").Value = 21749.25
$wsLTW.Range("M132").Value = -41542.571
$wsLTW.Range("N132").Value = -26809.25

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H9").Value = 11500
$wsWVR.Range("I9").Value = 20000
$wsWVR.Range("J9").Value = 3000
$wsWVR.Range("K9").Value = 20000
$wsWVR.Range("L9").Value = 3000
$wsWVR.Range("M9").Value = -19860

$wsWVR.Range("H14").Value = 6261.25
$wsWVR.Range("I14").Value = 6261.25
$wsWVR.Range("J14").Value = 0
$wsWVR.Range("K14").Value = 6261.25
$wsWVR.Range("L14").Value = 0
$wsWVR.Range("M14").Value = -6093.25

$wsWVR.Range("H122").Value = 1552.25
$wsWVR.Range("I122").Value = 1803.6
$wsWVR.Range("J122").Value = 1133.3334
$wsWVR.Range("K122").Value = 5410.799999999999
$wsWVR.Range("L122").Value = 3400.0002
$wsWVR.Range("M122").Value = -2960.799999999999
